$d = $word.ActiveDocument

$d.Content.Find.Execute("75-69=6", $true, $false, $false, $false, $false, $true, 1, $false, "19+69=88", 2) | Out-Null
$d.Content.Find.Execute("56-6=50", $true, $false, $false, $false, $false, $true, 1, $false, "29+22=51", 2) | Out-Null
$d.Content.Find.Execute("38-27=11", $true, $false, $false, $false, $false, $true, 1, $false, "62+12=74", 2) | Out-Null
$d.Content.Find.Execute("83-32=51", $true, $false, $false, $false, $false, $true, 1, $false, "32+5=37", 2) | Out-Null
$d.Content.Find.Execute("64-45=19", $true, $false, $false, $false, $false, $true, 1, $false, "16+42=58", 2) | Out-Null
$d.Content.Find.Execute("21+45=66", $true, $false, $false, $false, $false, $true, 1, $false, "49-41=8", 2) | Out-Null
$d.Content.Find.Execute("85-70=15", $true, $false, $false, $false, $false, $true, 1, $false, "9+90=99", 2) | Out-Null
$d.Content.Find.Execute("33+2=35", $true, $false, $false, $false, $false, $true, 1, $false, "33+45=78", 2) | Out-Null
$d.Content.Find.Execute("21+66=87", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=29", 2) | Out-Null
$d.Content.Find.Execute("15+15=30", $true, $false, $false, $false, $false, $true, 1, $false, "43-28=15", 2) | Out-Null
$d.Content.Find.Execute("41+24=65", $true, $false, $false, $false, $false, $true, 1, $false, "80+11=91", 2) | Out-Null
$d.Content.Find.Execute("36+45=81", $true, $false, $false, $false, $false, $true, 1, $false, "8+21=29", 2) | Out-Null
$d.Content.Find.Execute("0+74=74", $true, $false, $false, $false, $false, $true, 1, $false, "40+21=61", 2) | Out-Null
$d.Content.Find.Execute("33+36=69", $true, $false, $false, $false, $false, $true, 1, $false, "73-28=45", 2) | Out-Null
$d.Content.Find.Execute("74+10=84", $true, $false, $false, $false, $false, $true, 1, $false, "6+72=78", 2) | Out-Null
$d.Content.Find.Execute("33-6=27", $true, $false, $false, $false, $false, $true, 1, $false, "51+26=77", 2) | Out-Null
$d.Content.Find.Execute("61+9=70", $true, $false, $false, $false, $false, $true, 1, $false, "81-0=81", 2) | Out-Null
$d.Content.Find.Execute("15+44=59", $true, $false, $false, $false, $false, $true, 1, $false, "78-24=54", 2) | Out-Null
$d.Content.Find.Execute("89-24=65", $true, $false, $false, $false, $false, $true, 1, $false, "18+23=41", 2) | Out-Null
$d.Content.Find.Execute("88-9=79", $true, $false, $false, $false, $false, $true, 1, $false, "81-39=42", 2) | Out-Null
$d.Content.Find.Execute("59-39=20", $true, $false, $false, $false, $false, $true, 1, $false, "80+11=91", 2) | Out-Null
$d.Content.Find.Execute("52+27=79", $true, $false, $false, $false, $false, $true, 1, $false, "86-6=80", 2) | Out-Null
$d.Content.Find.Execute("9+82=91", $true, $false, $false, $false, $false, $true, 1, $false, "83-5=78", 2) | Out-Null
$d.Content.Find.Execute("22+0=22", $true, $false, $false, $false, $false, $true, 1, $false, "52+17=69", 2) | Out-Null
$d.Content.Find.Execute("60-55=5", $true, $false, $false, $false, $false, $true, 1, $false, "32-27=5", 2) | Out-Null
$d.Content.Find.Execute("44+46=90", $true, $false, $false, $false, $false, $true, 1, $false, "57+11=68", 2) | Out-Null
$d.Content.Find.Execute("91-14=77", $true, $false, $false, $false, $false, $true, 1, $false, "23-13=10", 2) | Out-Null
$d.Content.Find.Execute("70-22=48", $true, $false, $false, $false, $false, $true, 1, $false, "40+5=45", 2) | Out-Null
$d.Content.Find.Execute("41+53=94", $true, $false, $false, $false, $false, $true, 1, $false, "3+64=67", 2) | Out-Null
$d.Content.Find.Execute("46+33=79", $true, $false, $false, $false, $false, $true, 1, $false, "73-59=14", 2) | Out-Null
$d.Content.Find.Execute("8+56=64", $true, $false, $false, $false, $false, $true, 1, $false, "50+3=53", 2) | Out-Null
$d.Content.Find.Execute("57-16=41", $true, $false, $false, $false, $false, $true, 1, $false, "67+8=75", 2) | Out-Null
$d.Content.Find.Execute("0+73=73", $true, $false, $false, $false, $false, $true, 1, $false, "22+4=26", 2) | Out-Null
$d.Content.Find.Execute("41+31=72", $true, $false, $false, $false, $false, $true, 1, $false, "42+45=87", 2) | Out-Null
$d.Content.Find.Execute("74+2=76", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=6", 2) | Out-Null
$d.Content.Find.Execute("86-75=11", $true, $false, $false, $false, $false, $true, 1, $false, "60-6=54", 2) | Out-Null
$d.Content.Find.Execute("35-35=0", $true, $false, $false, $false, $false, $true, 1, $false, "53-38=15", 2) | Out-Null
$d.Content.Find.Execute("60-43=17", $true, $false, $false, $false, $false, $true, 1, $false, "17+43=60", 2) | Out-Null
$d.Content.Find.Execute("96-14=82", $true, $false, $false, $false, $false, $true, 1, $false, "49+4=53", 2) | Out-Null
$d.Content.Find.Execute("69-7=62", $true, $false, $false, $false, $false, $true, 1, $false, "99-98=1", 2) | Out-Null
$d.Content.Find.Execute("4+42=46", $true, $false, $false, $false, $false, $true, 1, $false, "69-53=16", 2) | Out-Null
$d.Content.Find.Execute("32-8=24", $true, $false, $false, $false, $false, $true, 1, $false, "30+40=70", 2) | Out-Null
$d.Content.Find.Execute("75-0=75", $true, $false, $false, $false, $false, $true, 1, $false, "63-56=7", 2) | Out-Null
$d.Content.Find.Execute("13+23=36", $true, $false, $false, $false, $false, $true, 1, $false, "16+62=78", 2) | Out-Null
$d.Content.Find.Execute("44-6=38", $true, $false, $false, $false, $false, $true, 1, $false, "91-10=81", 2) | Out-Null
$d.Content.Find.Execute("69-39=30", $true, $false, $false, $false, $false, $true, 1, $false, "61-18=43", 2) | Out-Null
$d.Content.Find.Execute("72-67=5", $true, $false, $false, $false, $false, $true, 1, $false, "98-27=71", 2) | Out-Null
$d.Content.Find.Execute("79-7=72", $true, $false, $false, $false, $false, $true, 1, $false, "8-1=7", 2) | Out-Null
$d.Content.Find.Execute("78-2=76", $true, $false, $false, $false, $false, $true, 1, $false, "37+34=71", 2) | Out-Null
$d.Content.Find.Execute("41+50=91", $true, $false, $false, $false, $false, $true, 1, $false, "9+13=22", 2) | Out-Null
$d.Content.Find.Execute("0+7=7", $true, $false, $false, $false, $false, $true, 1, $false, "98-2=96", 2) | Out-Null
$d.Content.Find.Execute("55+10=65", $true, $false, $false, $false, $false, $true, 1, $false, "68-30=38", 2) | Out-Null
$d.Content.Find.Execute("89-26=63", $true, $false, $false, $false, $false, $true, 1, $false, "36+44=80", 2) | Out-Null
$d.Content.Find.Execute("73+16=89", $true, $false, $false, $false, $false, $true, 1, $false, "72+6=78", 2) | Out-Null
$d.Content.Find.Execute("1+78=79", $true, $false, $false, $false, $false, $true, 1, $false, "22+12=34", 2) | Out-Null
$d.Content.Find.Execute("29+28=57", $true, $false, $false, $false, $false, $true, 1, $false, "39+48=87", 2) | Out-Null
$d.Content.Find.Execute("76-35=41", $true, $false, $false, $false, $false, $true, 1, $false, "48-1=47", 2) | Out-Null
$d.Content.Find.Execute("47+34=81", $true, $false, $false, $false, $false, $true, 1, $false, "26+22=48", 2) | Out-Null
$d.Content.Find.Execute("14+35=49", $true, $false, $false, $false, $false, $true, 1, $false, "24+49=73", 2) | Out-Null
$d.Content.Find.Execute("35-11=24", $true, $false, $false, $false, $false, $true, 1, $false, "72+27=99", 2) | Out-Null
$d.Content.Find.Execute("53+20=73", $true, $false, $false, $false, $false, $true, 1, $false, "82+8=90", 2) | Out-Null
$d.Content.Find.Execute("96-29=67", $true, $false, $false, $false, $false, $true, 1, $false, "75-72=3", 2) | Out-Null
$d.Content.Find.Execute("23-10=13", $true, $false, $false, $false, $false, $true, 1, $false, "65+6=71", 2) | Out-Null
$d.Content.Find.Execute("90+5=95", $true, $false, $false, $false, $false, $true, 1, $false, "36+22=58", 2) | Out-Null
$d.Content.Find.Execute("89-67=22", $true, $false, $false, $false, $false, $true, 1, $false, "24-6=18", 2) | Out-Null
$d.Content.Find.Execute("48+40=88", $true, $false, $false, $false, $false, $true, 1, $false, "39-6=33", 2) | Out-Null
$d.Content.Find.Execute("76+21=97", $true, $false, $false, $false, $false, $true, 1, $false, "7+9=16", 2) | Out-Null
$d.Content.Find.Execute("16-11=5", $true, $false, $false, $false, $false, $true, 1, $false, "42+38=80", 2) | Out-Null
$d.Content.Find.Execute("91-17=74", $true, $false, $false, $false, $false, $true, 1, $false, "12+1=13", 2) | Out-Null
$d.Content.Find.Execute("37-21=16", $true, $false, $false, $false, $false, $true, 1, $false, "8+85=93", 2) | Out-Null
$d.Content.Find.Execute("57+32=89", $true, $false, $false, $false, $false, $true, 1, $false, "21+38=59", 2) | Out-Null
$d.Content.Find.Execute("68-34=34", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=45", 2) | Out-Null
$d.Content.Find.Execute("1+49=50", $true, $false, $false, $false, $false, $true, 1, $false, "81-49=32", 2) | Out-Null
$d.Content.Find.Execute("20+49=69", $true, $false, $false, $false, $false, $true, 1, $false, "20+17=37", 2) | Out-Null
$d.Content.Find.Execute("77+16=93", $true, $false, $false, $false, $false, $true, 1, $false, "62-27=35", 2) | Out-Null
$d.Content.Find.Execute("63-39=24", $true, $false, $false, $false, $false, $true, 1, $false, "19+79=98", 2) | Out-Null
$d.Content.Find.Execute("11-2=9", $true, $false, $false, $false, $false, $true, 1, $false, "67-10=57", 2) | Out-Null
$d.Content.Find.Execute("6+20=26", $true, $false, $false, $false, $false, $true, 1, $false, "11+44=55", 2) | Out-Null
$d.Content.Find.Execute("38+22=60", $true, $false, $false, $false, $false, $true, 1, $false, "63-1=62", 2) | Out-Null
$d.Content.Find.Execute("29-24=5", $true, $false, $false, $false, $false, $true, 1, $false, "23+18=41", 2) | Out-Null
$d.Content.Find.Execute("57-39=18", $true, $false, $false, $false, $false, $true, 1, $false, "27-19=8", 2) | Out-Null
$d.Content.Find.Execute("51+25=76", $true, $false, $false, $false, $false, $true, 1, $false, "9+54=63", 2) | Out-Null
$d.Content.Find.Execute("98-57=41", $true, $false, $false, $false, $false, $true, 1, $false, "38-22=16", 2) | Out-Null
$d.Content.Find.Execute("97-28=69", $true, $false, $false, $false, $false, $true, 1, $false, "51-48=3", 2) | Out-Null
$d.Content.Find.Execute("57-46=11", $true, $false, $false, $false, $false, $true, 1, $false, "98-31=67", 2) | Out-Null
$d.Content.Find.Execute("70-65=5", $true, $false, $false, $false, $false, $true, 1, $false, "43+39=82", 2) | Out-Null
$d.Content.Find.Execute("76-69=7", $true, $false, $false, $false, $false, $true, 1, $false, "92-69=23", 2) | Out-Null
$d.Content.Find.Execute("88-75=13", $true, $false, $false, $false, $false, $true, 1, $false, "59-15=44", 2) | Out-Null
$d.Content.Find.Execute("9+34=43", $true, $false, $false, $false, $false, $true, 1, $false, "66-14=52", 2) | Out-Null
$d.Content.Find.Execute("20+77=97", $true, $false, $false, $false, $false, $true, 1, $false, "51+34=85", 2) | Out-Null
$d.Content.Find.Execute("20+2=22", $true, $false, $false, $false, $false, $true, 1, $false, "88-36=52", 2) | Out-Null
$d.Content.Find.Execute("1+26=27", $true, $false, $false, $false, $false, $true, 1, $false, "33+23=56", 2) | Out-Null
$d.Content.Find.Execute("84+8=92", $true, $false, $false, $false, $false, $true, 1, $false, "48-36=12", 2) | Out-Null
$d.Content.Find.Execute("90-68=22", $true, $false, $false, $false, $false, $true, 1, $false, "53-24=29", 2) | Out-Null
$d.Content.Find.Execute("59-53=6", $true, $false, $false, $false, $false, $true, 1, $false, "80-51=29", 2) | Out-Null
$d.Content.Find.Execute("17+37=54", $true, $false, $false, $false, $false, $true, 1, $false, "78-15=63", 2) | Out-Null
$d.Content.Find.Execute("82-60=22", $true, $false, $false, $false, $false, $true, 1, $false, "88-76=12", 2) | Out-Null
$d.Content.Find.Execute("24+50=74", $true, $false, $false, $false, $false, $true, 1, $false, "3+85=88", 2) | Out-Null
$d.Content.Find.Execute("87-61=26", $true, $false, $false, $false, $false, $true, 1, $false, "5+12=17", 2) | Out-Null
$d.Content.Find.Execute("55+21=76", $true, $false, $false, $false, $false, $true, 1, $false, "25+5=30", 2) | Out-Null
